$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "role"
$ws.Range("B1").Value = "content"

$ws.Range("D3").Select() | Out-Null
